# Apply the edit described in the commit: a new transaction row is
# inserted at row 3 of the sheet (pushing the existing rows 3..163 down
# to 4..164). The new row duplicates the formatting of row 2 (Excel's
# default behaviour when inserting a row above a formatted row) and
# contains a "Withdrawal / Credit Card / Tradeprof" transaction of
# 269.275 USD.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 3. Excel copies the formatting from the
# row above (row 2) onto the freshly inserted row, which matches the
# s="3"/s="4" styling seen on the new row in the target workbook.
$ws.Rows.Item(3).Insert()

# Populate the new row's data cells.
$ws.Range("E3").Value = "Withdrawal"
$ws.Range("N3").Value = "Credit Card"
$ws.Range("P3").Value = "Tradeprof"
$ws.Range("T3").Value = 269.27499999999998

# Best-effort: restore/update the view state (scroll position and
# selection) to match what was recorded after the edit.
$win = $excel.ActiveWindow
$win.ScrollRow = 154
$win.ScrollColumn = 11
$ws.Range("A165:XFD165").Select()
